$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists NBA players with their position and team. This edit
# reorders a subset of the data rows (the underlying rows were
# rearranged/re-sorted), so update each affected row's Name/Position/Team
# to reflect the new row order while leaving unaffected rows untouched.

$ws.Range("A7").Value = "Jalen Johnson"
$ws.Range("B7").Value = "PF"
$ws.Range("C7").Value = "Atlanta Hawks"

$ws.Range("A10").Value = "Jarrett Allen"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Cleveland Cavaliers"

$ws.Range("A11").Value = "Karl-Anthony Towns"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "New York Knicks"

$ws.Range("A13").Value = "Keegan Murray"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Sacramento Kings"

$ws.Range("A14").Value = "Daniel Gafford"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Dallas Mavericks"

$ws.Range("A15").Value = "Jalen Duren"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Detroit Pistons"

$ws.Range("A16").Value = "Taurean Prince"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Milwaukee Bucks"

$ws.Range("A17").Value = "Kevin Durant"
$ws.Range("B17").Value = "SF,PF"
$ws.Range("C17").Value = "Phoenix Suns"

$ws.Range("A18").Value = "Franz Wagner"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Orlando Magic"

$ws.Range("A19").Value = "Austin Reaves"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Los Angeles Lakers"
